$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.534.35"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "4.075.28"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.60"
$ws.Range("E5").Value = "  +5.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.02"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").Value = "4.069.38"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.697"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.771"
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.31"
$ws.Range("E12").Value = "  +15.01%  "
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.13"
$ws.Range("E14").Value = "  +5.07%  "
$ws.Range("D15").Value = "4.707.20"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "4.061.44"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.56"
$ws.Range("E17").Value = "  +4.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.89"
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "73.367.65"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "448.39"
$ws.Range("E22").Value = "  +5.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "98.50"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.47"
$ws.Range("E25").Value = "  +6.66%  "
$ws.Range("E26").Value = "  +3.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.31"
$ws.Range("E27").Value = "  +20.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.33"
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.96"
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.39"
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.88"
$ws.Range("E32").Value = "  +13.42%  "
$ws.Range("E33").Value = "  +4.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.65"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "687.86"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.72"
$ws.Range("E36").Value = "  +14.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "68.24"
$ws.Range("E37").Value = "  +4.61%  "
$ws.Range("D38").Value = "0.0₃0916"
$ws.Range("E38").Value = "  +11.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.450"
$ws.Range("E39").Value = "  +6.27%  "
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.24"
$ws.Range("E42").Value = "  +17.62%  "
$ws.Range("B43").Value = "Dai"
$ws.Range("C43").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.70"
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.58"
$ws.Range("E49").Value = "  +9.08%  "
$ws.Range("E50").Value = "  +5.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.32"
$ws.Range("E51").Value = "  +0.21%  "
